# Chapter 9 slides - fix typo (missing space) in the "some values" example
# textbox on slide 12 ("Multiple Bernoulli Document Representation"):
#   "...spam)=1,P(dinner..."  ->  "...spam)=1, P(dinner..."
#
# Re-setting the whole two-paragraph string (joined with a CR so it keeps
# its original two <a:p> paragraphs) reproduces the corrected wording, and
# PowerPoint's own autofit ("shrink/grow shape to fit text" textbox) then
# needs its stored width nudged back up to the size PowerPoint itself would
# have relaid the (now very slightly wider) text out to - the host already
# recomputes the autofit height on save, but not the autofit width.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(12)
$shp = $s.Shapes.Item(3)

$shp.TextFrame.TextRange.Text = "some values: P(spam)=3/10, P(not spam)=7/10, P(the|spam)=1, " + "`r" + "P(the|not spam)=1, P(dinner|spam)=0, P(dinner|not spam)=1/7,…"

# Nudge the autofit ("Resize shape to fit text") width back out to match
# the slightly-longer corrected line (height is relaid out automatically).
$shp.Width = 542.4282
